# Update bee action function - mark rows as "Merged to master"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# Row 10 (LongThrowerAnt): status changed from "Pull request open" to "Merged to master"
$ws.Range("E10").Value = "Merged to master"

# Row 12 (ThrowerAnt_Test): status set to "Merged to master" (was empty)
$ws.Range("E12").Value = "Merged to master"

# Row 13 (LongThrowerAnt_Test): status set to "Merged to master" (was empty)
$ws.Range("E13").Value = "Merged to master"

# Update the active selection / view to E14, scrolled back to top
$ws.Activate()
$ws.Range("E14").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
